# Implements unmodified logit for transportation allocations: replace the
# formulas that pulled logit exponents from the
# 'A54.tranSubsector_logit_revised' sheet with hardcoded literal values on
# the TTLE sheet.

$wb = $excel.ActiveWorkbook

$wsTTLE = $wb.Worksheets.Item("TTLE")

# Column B and C, rows 2-7: replace formulas with plain numeric values.
$wsTTLE.Range("B2").Value = -20
$wsTTLE.Range("C2").Value = -20

$wsTTLE.Range("B3").Value = -20
$wsTTLE.Range("C3").Value = -20

$wsTTLE.Range("B4").Value = -0.02
$wsTTLE.Range("C4").Value = -0.02

$wsTTLE.Range("B5").Value = -0.02
$wsTTLE.Range("C5").Value = -0.02

$wsTTLE.Range("B6").Value = -0.02
$wsTTLE.Range("C6").Value = -0.02

$wsTTLE.Range("B7").Value = -20
$wsTTLE.Range("C7").Value = -20

# Update sheet selections / active sheet to match the reviewed workbook
# state: "About" sheet selection moves to C21 and is no longer the
# tab-selected sheet; "TTLE" becomes the tab-selected / active sheet with
# its selection spanning B2:C7.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Range("C21").Select() | Out-Null

$wsTTLE.Activate() | Out-Null
$wsTTLE.Range("B2:C7").Select() | Out-Null
